$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 139.8
$ws.Range("J12").Value = 139.5
$ws.Range("L12").Value = 139.5
$ws.Range("N12").Value = -479.5
$ws.Range("H53").Value = 1016.5
$ws.Range("I53").Value = 906.5333000000001
$ws.Range("K53").Value = 906.5333000000001
$ws.Range("M53").Value = -269.5333000000001
$ws.Range("H64").Value = 3422.9285
$ws.Range("I64").Value = 3330.3333
$ws.Range("K64").Value = 3330.3333
$ws.Range("M64").Value = -3082.3333
$ws.Range("H67").Value = 3422.9285
$ws.Range("I67").Value = 3330.3333
$ws.Range("K67").Value = 3330.3333
$ws.Range("M67").Value = -2472.3333
$ws.Range("H76").Value = 3062.6667
$ws.Range("I76").Value = 3062.6667
$ws.Range("K76").Value = 3062.6667
$ws.Range("M76").Value = -2747.6667
$ws.Range("H79").Value = 3062.6667
$ws.Range("I79").Value = 3062.6667
$ws.Range("K79").Value = 3062.6667
$ws.Range("M79").Value = -1970.6667
$ws.Range("H86").Value = 3729.6155
$ws.Range("I86").Value = 3277
$ws.Range("J86").Value = 4257.6665
$ws.Range("K86").Value = 3277
$ws.Range("L86").Value = 4257.6665
$ws.Range("M86").Value = -2154
$ws.Range("N86").Value = -6503.6665
$ws.Range("H89").Value = 3729.6155
$ws.Range("I89").Value = 3277
$ws.Range("J89").Value = 4257.6665
$ws.Range("K89").Value = 16385
$ws.Range("L89").Value = 21288.3325
$ws.Range("M89").Value = -10769
$ws.Range("N89").Value = -32520.3325
$ws.Range("H98").Value = 2223
$ws.Range("I98").Value = 2286
$ws.Range("J98").Value = 2097
$ws.Range("K98").Value = 2286
$ws.Range("L98").Value = 2097
$ws.Range("M98").Value = -788
$ws.Range("N98").Value = -5093
$ws.Range("H121").Value = 1214.2858
$ws.Range("J121").Value = 1214.2858
$ws.Range("L121").Value = 3642.8574
$ws.Range("N121").Value = -7136.857400000001
$ws.Range("H122").Value = 2223
$ws.Range("I122").Value = 2286
$ws.Range("J122").Value = 2097
$ws.Range("K122").Value = 6858
$ws.Range("L122").Value = 6291
$ws.Range("M122").Value = -4408
$ws.Range("N122").Value = -11191
$ws.Range("H137").Value = 2079.7058
$ws.Range("I137").Value = 1180.875
$ws.Range("K137").Value = 3542.625
$ws.Range("M137").Value = -992.625
$ws.Range("H138").Value = 2155.72
$ws.Range("J138").Value = 2235.7527
$ws.Range("L138").Value = 6707.2581
$ws.Range("N138").Value = -16987.2581
$ws.Range("H141").Value = 11856.9
$ws.Range("I141").Value = 12729.889
$ws.Range("K141").Value = 38189.667
$ws.Range("M141").Value = -33009.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1261.3684
$ws.Range("I2").Value = 716.7778
$ws.Range("J2").Value = 1751.5
$ws.Range("K2").Value = 716.7778
$ws.Range("L2").Value = 1751.5
$ws.Range("M2").Value = -603.7778
$ws.Range("N2").Value = -1977.5
$ws.Range("H32").Value = 2654.0146
$ws.Range("I32").Value = 2726.5173
$ws.Range("J32").Value = 2233.5
$ws.Range("K32").Value = 2726.5173
$ws.Range("L32").Value = 2233.5
$ws.Range("M32").Value = -2439.5173
$ws.Range("N32").Value = -2807.5
$ws.Range("H61").Value = 964.7931
$ws.Range("I61").Value = 619.4091
$ws.Range("J61").Value = 2050.2856
$ws.Range("K61").Value = 619.4091
$ws.Range("L61").Value = 2050.2856
$ws.Range("M61").Value = -407.4091
$ws.Range("N61").Value = -2474.2856
$ws.Range("H116").Value = 1261.3684
$ws.Range("I116").Value = 716.7778
$ws.Range("J116").Value = 1751.5
$ws.Range("K116").Value = 716.7778
$ws.Range("L116").Value = 1751.5
$ws.Range("M116").Value = 1577.2222
$ws.Range("N116").Value = -6339.5
$ws.Range("H123").Value = 74607.25
$ws.Range("J123").Value = 74607.25
$ws.Range("L123").Value = 74607.25
$ws.Range("N123").Value = -84407.25
$ws.Range("H136").Value = 964.7931
$ws.Range("I136").Value = 619.4091
$ws.Range("J136").Value = 2050.2856
$ws.Range("K136").Value = 1858.2273
$ws.Range("L136").Value = 6150.8568
$ws.Range("M136").Value = 691.7727
$ws.Range("N136").Value = -11250.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1261.3684
$ws.Range("I3").Value = 716.7778
$ws.Range("J3").Value = 1751.5
$ws.Range("K3").Value = 716.7778
$ws.Range("L3").Value = 1751.5
$ws.Range("M3").Value = -602.7778
$ws.Range("N3").Value = -1979.5
$ws.Range("H134").Value = 9596.808000000001
$ws.Range("I134").Value = 6492.909
$ws.Range("K134").Value = 19478.727
$ws.Range("M134").Value = -16943.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1355.3914
$ws.Range("J31").Value = 1723.1
$ws.Range("L31").Value = 1723.1
$ws.Range("N31").Value = -2313.1
$ws.Range("H34").Value = 1355.3914
$ws.Range("J34").Value = 1723.1
$ws.Range("L34").Value = 1723.1
$ws.Range("N34").Value = -2127.1
$ws.Range("H62").Value = 8698578
$ws.Range("I62").Value = 3065
$ws.Range("J62").Value = 66668668
$ws.Range("K62").Value = 3065
$ws.Range("L62").Value = 66668668
$ws.Range("M62").Value = -2441
$ws.Range("N62").Value = -66669916
$ws.Range("H65").Value = 8698578
$ws.Range("I65").Value = 3065
$ws.Range("J65").Value = 66668668
$ws.Range("K65").Value = 15325
$ws.Range("L65").Value = 333343340
$ws.Range("M65").Value = -12205
$ws.Range("N65").Value = -333349580
$ws.Range("H134").Value = 10102418
$ws.Range("I134").Value = 12821869
$ws.Range("K134").Value = 38465607
$ws.Range("M134").Value = -38463072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6197.1816
$ws.Range("I3").Value = 4816.9
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 14450.7
$ws.Range("L3").Value = 60000
$ws.Range("M3").Value = -14338.7
$ws.Range("N3").Value = -60224
$ws.Range("H23").Value = 182.78572
$ws.Range("I23").Value = 97
$ws.Range("K23").Value = 291
$ws.Range("M23").Value = -56
$ws.Range("H41").Value = 198
$ws.Range("I41").Value = 198
$ws.Range("K41").Value = 594
$ws.Range("M41").Value = -256
$ws.Range("H121").Value = 772.3684
$ws.Range("I121").Value = 277.33334
$ws.Range("J121").Value = 865.1875
$ws.Range("K121").Value = 832.0000200000001
$ws.Range("L121").Value = 2595.5625
$ws.Range("M121").Value = 477.9999799999999
$ws.Range("N121").Value = -5215.5625
$ws.Range("H141").Value = 4378
$ws.Range("I141").Value = 4222.5
$ws.Range("K141").Value = 12667.5
$ws.Range("M141").Value = -7487.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2627.0952
$ws.Range("I122").Value = 1714.3077
$ws.Range("K122").Value = 5142.9231
$ws.Range("M122").Value = -2692.9231
$ws.Range("H132").Value = 2509.96
$ws.Range("I132").Value = 2028.1875
$ws.Range("K132").Value = 6084.5625
$ws.Range("M132").Value = -3554.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3021.2222
$ws.Range("I40").Value = 2740.8572
$ws.Range("K40").Value = 2740.8572
$ws.Range("M40").Value = -2604.8572
$ws.Range("H61").Value = 2092.5264
$ws.Range("I61").Value = 1473.3077
$ws.Range("J61").Value = 3434.1667
$ws.Range("K61").Value = 1473.3077
$ws.Range("L61").Value = 3434.1667
$ws.Range("M61").Value = -1271.3077
$ws.Range("N61").Value = -3838.1667
$ws.Range("H68").Value = 1345
$ws.Range("I68").Value = 1096.8
$ws.Range("J68").Value = 1699.5714
$ws.Range("K68").Value = 1096.8
$ws.Range("L68").Value = 1699.5714
$ws.Range("M68").Value = -347.8
$ws.Range("N68").Value = -3197.5714
$ws.Range("H71").Value = 1345
$ws.Range("I71").Value = 1096.8
$ws.Range("J71").Value = 1699.5714
$ws.Range("K71").Value = 5484
$ws.Range("L71").Value = 8497.857
$ws.Range("M71").Value = -1740
$ws.Range("N71").Value = -15985.857
$ws.Range("H93").Value = 1424.75
$ws.Range("I93").Value = 1281.5454
$ws.Range("K93").Value = 1281.5454
$ws.Range("M93").Value = -33.54539999999997
$ws.Range("H113").Value = 2092.5264
$ws.Range("I113").Value = 1473.3077
$ws.Range("J113").Value = 3434.1667
$ws.Range("K113").Value = 1473.3077
$ws.Range("L113").Value = 3434.1667
$ws.Range("M113").Value = 696.6922999999999
$ws.Range("N113").Value = -7774.1667
$ws.Range("H116").Value = 25680
$ws.Range("J116").Value = 25680
$ws.Range("L116").Value = 25680
$ws.Range("N116").Value = -34858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 949.125
$ws.Range("I113").Value = 398.25
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1194.75
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = 975.25
$ws.Range("N113").Value = -8840
